# Update cfb_weather.xlsx with Timestamp 2025-12-05T10:01:22.193876
# Re-applies refreshed weather-pull data for each scheduled game on both
# the "FBS" and "Other" sheets (wind direction/speed, temps, forecast
# deltas) plus the new pull Timestamp.

$wb = $excel.ActiveWorkbook

# ---- Sheet "FBS" ----
$ws1 = $wb.Worksheets.Item("FBS")

# Row 2 - UNLV @ Boise State
$ws1.Range("M2").Value  = "W"
$ws1.Range("N2").Value  = "WNW"
$ws1.Range("O2").Value  = 44.6
$ws1.Range("P2").Value  = 4.3
$ws1.Range("Q2").Value  = "WNW"
$ws1.Range("R2").Value  = 1.3
$ws1.Range("S2").Value  = -1.5
$ws1.Range("U2").Value  = -2.5
$ws1.Range("Z2").Value  = -105
$ws1.Range("AK2").Value = "2025-12-05T10:01:22.193876"

# Row 3 - Troy @ James Madison
$ws1.Range("N3").Value  = "NNE"
$ws1.Range("O3").Value  = 22.22
$ws1.Range("P3").Value  = 3
$ws1.Range("Q3").Value  = "WNW"
$ws1.Range("S3").Value  = -0.97
$ws1.Range("T3").Value  = -1.22
$ws1.Range("U3").Value  = -3.2
$ws1.Range("AK3").Value = "2025-12-05T10:01:22.193876"

# Row 4 - Kennesaw State @ Jacksonville State
$ws1.Range("O4").Value  = 43.88
$ws1.Range("P4").Value  = 4.7
$ws1.Range("U4").Value  = -0.1
$ws1.Range("AK4").Value = "2025-12-05T10:01:22.193876"

# Row 5 - North Texas @ Tulane
$ws1.Range("M5").Value  = "ENE"
$ws1.Range("N5").Value  = "NNW"
$ws1.Range("O5").Value  = 53
$ws1.Range("P5").Value  = 1.7
$ws1.Range("Q5").Value  = "NNW"
$ws1.Range("R5").Value  = 0
$ws1.Range("U5").Value  = -8.9
$ws1.Range("AK5").Value = "2025-12-05T10:01:22.193876"

# Row 6 - Miami (OH) @ Western Michigan
$ws1.Range("O6").Value  = 32.35999999999999
$ws1.Range("P6").Value  = 6.5
$ws1.Range("U6").Value  = -5.5
$ws1.Range("AK6").Value = "2025-12-05T10:01:22.193876"

# Row 7 - Duke @ Virginia
$ws1.Range("O7").Value  = 32.48
$ws1.Range("P7").Value  = 3.8
$ws1.Range("U7").Value  = -0.8
$ws1.Range("AK7").Value = "2025-12-05T10:01:22.193876"

# ---- Sheet "Other" ----
$ws2 = $wb.Worksheets.Item("Other")

# Row 2 - Yale vs Montana State
$ws2.Range("Q2").Value = 38.66
$ws2.Range("R2").Value = 17.9
$ws2.Range("U2").Value = -6.5

# Row 3 - South Dakota vs Mercer
$ws2.Range("O3").Value = "NE"
$ws2.Range("P3").Value = "ENE"
$ws2.Range("Q3").Value = 50.59999999999999
$ws2.Range("R3").Value = 5.5
$ws2.Range("S3").Value = "ENE"

# Row 4 - Villanova vs Lehigh
$ws2.Range("O4").Value = "ENE"
$ws2.Range("Q4").Value = 39.2
$ws2.Range("R4").Value = 6.3

# Row 5 - North Dakota vs Tarleton State
$ws2.Range("P5").Value = "NE"
$ws2.Range("Q5").Value = 59.84
$ws2.Range("R5").Value = 10

# Row 6 - Abilene Christian vs Stephen F. Austin
$ws2.Range("P6").Value = "NNW"
$ws2.Range("Q6").Value = 54.38
$ws2.Range("R6").Value = 6.8

# Row 7 - South Dakota State vs Montana
$ws2.Range("O7").Value = "NE"
$ws2.Range("Q7").Value = 39.2
$ws2.Range("R7").Value = 9.4
$ws2.Range("S7").Value = "NE"
$ws2.Range("T7").Value = 0.4

# Row 8 - Rhode Island vs UC Davis
$ws2.Range("O8").Value = "W"
$ws2.Range("P8").Value = "SSW"
$ws2.Range("Q8").Value = 42.62
$ws2.Range("R8").Value = 2.3
$ws2.Range("S8").Value = "W"
